# Natmi following Dr Hou advice
# Recomputed Gdf9-Bmpr1a ligand/receptor statistics for the three sending
# clusters (ECs, FAPs, sCs) against the three target clusters, adding the
# previously-missing "sCs" sending-cluster rows (8-10) and refreshing the
# values for the existing rows (2-7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r2 = @("ECs", "Gdf9", "Bmpr1a", "ECs", 3, 1, 0.9029553333333333, 2.708866, 0.2353208431348719, 0.2353208431348719, 3, 1, 9.918859, 29.756577, 0.160764128269069, 0.160764128269069, 8.956286634631333, 80.606579711682, 0.03783115021012002, 0.03783115021012003)
$r3 = @("ECs", "Gdf9", "Bmpr1a", "FAPs", 3, 1, 0.9029553333333333, 2.708866, 0.2353208431348719, 0.2353208431348719, 3, 1, 29.20351433333333, 87.610543, 0.47332838627826, 0.4733283862782601, 26.36946901935978, 237.325221174238, 0.1113840349386685, 0.1113840349386685)
$r4 = @("ECs", "Gdf9", "Bmpr1a", "sCs", 3, 1, 0.9029553333333333, 2.708866, 0.2353208431348719, 0.2353208431348719, 3, 1, 22.575837, 67.72751099999999, 0.3659074854526709, 0.3659074854526709, 20.384972423614, 183.464751812526, 0.0861056579860834, 0.08610565798608341)
$r5 = @("FAPs", "Gdf9", "Bmpr1a", "ECs", 3, 1, 1.446898, 4.340694, 0.3770787376970584, 0.3770787376970584, 3, 1, 9.918859, 29.756577, 0.160764128269069, 0.160764128269069, 14.351577249382, 129.164195244438, 0.06062073455466853, 0.06062073455466853)
$r6 = @("FAPs", "Gdf9", "Bmpr1a", "FAPs", 3, 1, 1.446898, 4.340694, 0.3770787376970584, 0.3770787376970584, 3, 1, 29.20351433333333, 87.610543, 0.47332838627826, 0.4733283862782601, 42.25450648187134, 380.290558336842, 0.1784820704139919, 0.178482070413992)
$r7 = @("FAPs", "Gdf9", "Bmpr1a", "sCs", 3, 1, 1.446898, 4.340694, 0.3770787376970584, 0.3770787376970584, 3, 1, 22.575837, 67.72751099999999, 0.3659074854526709, 0.3659074854526709, 32.66493340362599, 293.984400632634, 0.1379759327283979, 0.1379759327283979)
$r8 = @("sCs", "Gdf9", "Bmpr1a", "ECs", 2, 0.6666666666666666, 1.487271, 4.461812999999999, 0.3876004191680696, 0.3876004191680697, 3, 1, 9.918859, 29.756577, 0.160764128269069, 0.160764128269069, 14.752031343789, 132.768282094101, 0.06231224350428047, 0.06231224350428047)
$r9 = @("sCs", "Gdf9", "Bmpr1a", "FAPs", 2, 0.6666666666666666, 1.487271, 4.461812999999999, 0.3876004191680696, 0.3876004191680697, 3, 1, 29.20351433333333, 87.610543, 0.47332838627826, 0.4733283862782601, 43.43353996605099, 390.9018596944589, 0.1834622809255996, 0.1834622809255996)
$r10 = @("sCs", "Gdf9", "Bmpr1a", "sCs", 2, 0.6666666666666666, 1.487271, 4.461812999999999, 0.3876004191680696, 0.3876004191680697, 3, 1, 22.575837, 67.72751099999999, 0.3659074854526709, 0.3659074854526709, 33.57638767082699, 302.1874890374429, 0.1418258947381896, 0.1418258947381896)

$data = @($r2, $r3, $r4, $r5, $r6, $r7, $r8, $r9, $r10)

$rowIndex = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $row[$c]
    }
    $rowIndex++
}